$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Lista de Tabelas"
$ws2 = $wb.Worksheets.Item(2)   # "IHH - Massa de rendimentos"

# --- Update the year reference 2014 -> 2015 in the caption text used on both sheets ---
$ws1.Range("A2").Value = " Tabela 1 - Evolução do Índice de Concentração IHH do Total de Rendimentos das Atividades Relacionadas por UF entre 2007 e 2015"
$ws2.Range("A1").Value = " Tabela 1 - Evolução do Índice de Concentração IHH do Total de Rendimentos das Atividades Relacionadas por UF entre 2007 e 2015"

# --- Shift the year header row on sheet2: drop 2010, extend through 2015 ---
$ws2.Range("E2").Value = 2011
$ws2.Range("F2").Value = 2012
$ws2.Range("G2").Value = 2013
$ws2.Range("H2").Value = 2014
$ws2.Range("I2").Value = 2015

# --- Update the IHH data row (B3:I3) with refreshed figures ---
$ws2.Range("B3").Value = 0.18384039004020575
$ws2.Range("C3").Value = 0.18013115688323908
$ws2.Range("D3").Value = 0.1647299461193342
$ws2.Range("E3").Value = 0.18345567381763003
$ws2.Range("F3").Value = 0.17255709871391384
$ws2.Range("G3").Value = 0.16852289273509774
$ws2.Range("H3").Value = 0.17838328167136369
# I3 (0.17042321055268089) is unchanged

# --- Move the selection cursor on sheet1, then make sheet2 the active/selected tab ---
$ws1.Range("A3").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("H9").Select() | Out-Null
